# Generate Report for Handback
# Updates timestamps/priority produced by a re-run of the handback report
# generation (shared strings shown in the Overview / zh-cn / de-de tables).

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for e318f462-... (row 4) and
# f49b2460-... (row 5) both held the same timestamp string, so both
# cells move together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-21 16:15:57"
$wsOverview.Range("G5").Value = "2016-08-21 16:15:57"

# --- zh-cn sheet ------------------------------------------------------
# Priority ("ht" -> "mt") and the handoff/handback datetimes for the
# e318f462-... (row 4) and f49b2460-... (row 5) rows, which shared the
# same values.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-21 16:15:53"
$wsZhCn.Range("K4").Value = "2016-08-21 16:16:16"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H5").Value = "2016-08-21 16:15:53"
$wsZhCn.Range("K5").Value = "2016-08-21 16:16:16"

# --- de-de sheet --------------------------------------------------
# Handoff datetime (shared with the Overview date above) and handback
# datetime for the e318f462-... (row 4) and f49b2460-... (row 5) rows.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-21 16:15:57"
$wsDeDe.Range("K4").Value = "2016-08-21 16:16:22"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H5").Value = "2016-08-21 16:15:57"
$wsDeDe.Range("K5").Value = "2016-08-21 16:16:22"
